# TC29_Canine_Filter_Breed-Gordon.xlsx -- "Fixed ICDC breed all testcases"
#
# Replace the old StatQuery (column C, rows 2-4) with the new query that
# returns Programs / Studies / Cases / Samples / Case Files / Study Files
# counts, and update the sheet view (zoom + selection) to match the
# author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = @"
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Gordon Setter']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``
"@
$newStatQuery = $newStatQuery.TrimEnd("`r", "`n")

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the saved view state: zoom 55 -> 85, scroll so row 3 is at the
# top, and move the selection to B4.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
